# The edit moves the auto "_GoBack" bookmark (last-edit marker) from the
# end of the document up into the "Opearation" bullet, where the author's
# cursor last made an edit: the run "Opearation: Administrator can turn
# off connection, ..." is split after "Ope" (dropping the stray "a" so it
# reads "Ope" + "ration: ..." = "Operation: ..."), with an empty
# _GoBack bookmark inserted at the split point. The bookmark that used to
# sit at the very end of the document (after "So we have to design status
# tab item for each slot!") is removed as part of the move.

$d = $word.ActiveDocument

# Locate the target run via Find.
$found = $d.Content.Duplicate
$found.Find.Execute("Opearation", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)

# Split point: after "Ope" (3 chars) into the matched word.
$splitPoint = $found.Start + 3

# Drop the extra "a" right after "Ope" so the text reads
# "Ope" | "ration: ..." (i.e. "Operation: ...").
$stray = $d.Range($splitPoint, $splitPoint + 1)
$stray.Text = ""

# Re-home the (unique, auto-maintained) "_GoBack" bookmark here; Word
# enforces unique bookmark names, so adding it here removes/moves it from
# wherever it previously lived (the end of the document).
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
